$d = $word.ActiveDocument

$pairs = @(
    @("23÷5=4, 3", "58÷5=11, 3"),
    @("21÷5=4, 1", "17÷3=5, 2"),
    @("27÷9=3, 0", "79÷4=19, 3"),
    @("37÷8=4, 5", "77÷3=25, 2"),
    @("30÷6=5, 0", "72÷3=24, 0"),
    @("64÷8=8, 0", "92÷2=46, 0"),
    @("22÷5=4, 2", "80÷8=10, 0"),
    @("44÷3=14, 2", "85÷3=28, 1"),
    @("51÷6=8, 3", "16÷8=2, 0"),
    @("89÷8=11, 1", "84÷2=42, 0"),
    @("75÷2=37, 1", "15÷8=1, 7"),
    @("47÷3=15, 2", "40÷6=6, 4"),
    @("72÷8=9, 0", "66÷4=16, 2"),
    @("42÷5=8, 2", "54÷7=7, 5"),
    @("70÷5=14, 0", "38÷2=19, 0"),
    @("19÷6=3, 1", "57÷7=8, 1"),
    @("97÷5=19, 2", "91÷7=13, 0"),
    @("61÷3=20, 1", "21÷8=2, 5"),
    @("70÷9=7, 7", "17÷3=5, 2"),
    @("83÷7=11, 6", "66÷8=8, 2"),
    @("70÷7=10, 0", "87÷9=9, 6"),
    @("99÷5=19, 4", "84÷8=10, 4"),
    @("61÷7=8, 5", "51÷7=7, 2"),
    @("66÷9=7, 3", "87÷4=21, 3"),
    @("99÷4=24, 3", "14÷9=1, 5")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false,
                         $true, 1, $false, $new, 2)
}
